{"js": "// Change 1 & 2: the arraignment date and the \"paid in full by\" date both\n// move from December 15, 2021 -> December 17, 2021. Both instances of the\n// literal string are updated in place so each run keeps its own formatting.\nconst dateResults = context.document.body.search(\"December 15, 2021\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"December 17, 2021\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Change 3: in the \"Proof of Financial Responsibility\" paragraph, the old\n// sentence \"The Defendant showed the Court proof of responsibility during\n// the proceeding.\" is replaced with the new, longer sentence. The diff also\n// shows the lone space run that used to sit between \"Responsibility.  \" and\n// that sentence becoming empty (i.e. the extra space is dropped), so the\n// final paragraph text goes from three spaces to two spaces before \"The\".\nconst oldSentence =\n  \"The Defendant showed the Court proof of responsibility during the proceeding.\";\nconst newSentence =\n  \"The Defendant did not show proof of financial responsibility at the time of the offense \" +\n  \"or during the proceeding, but may show proof to Clerk of Court at any time prior to the \" +\n  \"submission of this matter to the Ohio Bureau of Motor Vehicles.\";\n\n// Locate the lone-space run right before the sentence, and the sentence\n// itself, so we can drop the stray space and swap the sentence text while\n// each keeps its own run formatting.\nconst headResults = context.document.body.search(\"Responsibility.  \", { matchCase: true });\nheadResults.load(\"text\");\nawait context.sync();\n\nconst sentResults = context.document.body.search(oldSentence, { matchCase: true });\nsentResults.load(\"text\");\nawait context.sync();\n\nif (headResults.items.length > 0 && sentResults.items.length > 0) {\n  const afterHead = headResults.items[0].getRange(\"After\");\n  const sentStart = sentResults.items[0].getRange(\"Start\");\n  const spaceRange = afterHead.expandTo(sentStart);\n  spaceRange.insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Re-search for the sentence (its range may have shifted after the space\n// removal) and swap in the new wording.\nconst sentResults2 = context.document.body.search(oldSentence, { matchCase: true });\nsentResults2.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < sentResults2.items.length; i++) {\n  sentResults2.items[i].insertText(newSentence, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Changes 1 & 2: both the arraignment date and the \"paid in full by\" date\n# move from December 15, 2021 -> December 17, 2021.\n$find = $d.Content.Find\n$find.Text = \"December 15, 2021\"\n$find.Replacement.Text = \"December 17, 2021\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# Change 3: in the \"Proof of Financial Responsibility\" paragraph, the\n# sentence \"The Defendant showed the Court proof of responsibility during\n# the proceeding.\" is replaced with a longer sentence, and the lone extra\n# space that used to sit right before that sentence is dropped (so the\n# paragraph goes from three spaces to two spaces before \"The\").\n\n# Step 1: find the end of \"Responsibility.  \" and drop the single space\n# character that immediately follows it.\n$headRange = $d.Content.Duplicate\n$headFound = $headRange.Find.Execute(\"Responsibility.  \")\nif ($headFound) {\n    $spaceRange = $d.Range($headRange.End, $headRange.End + 1)\n    if ($spaceRange.Text -eq \" \") {\n        $spaceRange.Text = \"\"\n    }\n}\n\n# Step 2: re-find the old sentence (now correctly positioned after the\n# space was removed) and swap in the new wording, keeping that run's own\n# formatting.\n$sentRange = $d.Content.Duplicate\n$sentFound = $sentRange.Find.Execute(\"The Defendant showed the Court proof of responsibility during the proceeding.\")\nif ($sentFound) {\n    $sentRange.Text = \"The Defendant did not show proof of financial responsibility at the time of the offense or during the proceeding, but may show proof to Clerk of Court at any time prior to the submission of this matter to the Ohio Bureau of Motor Vehicles.\"\n}\n"}
